$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws.Range("B3").Value2 = "6.0.0"

# Date updated
$ws.Range("B8").Value2 = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$ws.Range("B9").Value2 = "Alvearie Team"

# Remove the duplicate "Contact" / "No display for ContactDetail" row (old row 11)
$ws.Rows("11:11").Delete()

# Old row 10 (still row 10 after the delete) becomes Jurisdiction / United States of America
$ws.Range("A10").Value2 = "Jurisdiction"
$ws.Range("B10").Value2 = "United States of America"

# Case Sensitive value set to text "true" (row 15 shifted up to row 14 after the delete).
# A direct Value2 = "true" assignment gets auto-coerced to a Boolean, so stage the
# text in a scratch cell via a formula and paste the *value* back in, which keeps it
# as a genuine text string instead of a Boolean.
$ws.Range("Z1").Formula = "=""true"""
$ws.Range("Z1").Copy()
$ws.Range("B14").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
